# Update gh-pages to output generated at 456a3b4
# Apply the updated "想去人数" (interest count) figures to column F
# on both the "展览" (exhibition) sheet and the "全部类型" (all types)
# sheet, which carry the same underlying data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8312
    3  = 7736
    9  = 115
    12 = 701
    13 = 125
    14 = 1297
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
